# Insert a new row before row 232, shifting existing rows 232:268 down to 233:269.
# The new row 232 duplicates the data that was (prior to the insert) in row 232,
# except for the date in column D which becomes the new date 45015 (2023-03-30).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows(232).Insert()

$ws.Cells.Item(232, 1).Value = 5
$ws.Cells.Item(232, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(232, 3).Value = "Maule"
$ws.Cells.Item(232, 4).Value = 45015
$ws.Cells.Item(232, 5).Value = 7
$ws.Cells.Item(232, 6).Value = 100112017
$ws.Cells.Item(232, 7).Value = "Apio"
$ws.Cells.Item(232, 8).Value = "Americana (o)"
$ws.Cells.Item(232, 9).Value = "Primera"
$ws.Cells.Item(232, 10).Value = 500
$ws.Cells.Item(232, 11).Value = 8000
$ws.Cells.Item(232, 12).Value = 8000
$ws.Cells.Item(232, 13).Value = 8000
$ws.Cells.Item(232, 14).Value = "`$/docena de matas"
$ws.Cells.Item(232, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(232, 16).Value = 1333
$ws.Cells.Item(232, 17).Value = 6
$ws.Cells.Item(232, 18).Value = "Hortaliza"

$ws.Cells.Item(232, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
